$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while guaranteeing it stays text-typed
# (protects numeric-looking strings like "517.99" or "1.00" from being
# auto-coerced to a Double by the COM Value setter) and without leaving the
# cells style/number-format changed from its original General/no-style state.
function Set-TextValue($range, [string]$text) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

$ws.Range("D2").Value = "58.192.29"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.599.77"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "517.99"
$ws.Range("E5").Value = "  -0.43%  "
Set-TextValue $ws.Range("D6") "143.05"
$ws.Range("E6").Value = "  -0.43%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue $ws.Range("D8") "0.566"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "2.599.90"
$ws.Range("E9").Value = "  -1.67%  "
Set-TextValue $ws.Range("D10") "6.86"
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "3.056.73"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "58.160.50"
$ws.Range("E15").Value = "  -0.36%  "
Set-TextValue $ws.Range("D16") "20.38"
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "2.565.02"
$ws.Range("E18").Value = "  -3.03%  "
Set-TextValue $ws.Range("D19") "342.97"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("E20").Value = "  -1.86%  "
Set-TextValue $ws.Range("D21") "10.30"
$ws.Range("E21").Value = "  -1.41%  "
Set-TextValue $ws.Range("D22") "6.35"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("E23").Value = "  +0.31%  "
Set-TextValue $ws.Range("D24") "66.27"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "2.701.46"
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").Value = "0.0₃0746"
$ws.Range("E30").Value = "  -6.02%  "
Set-TextValue $ws.Range("D31") "0.999"
$ws.Range("E31").Value = "  +0.02%  "
Set-TextValue $ws.Range("D32") "6.23"
$ws.Range("E32").Value = "  -5.42%  "
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  -0.16%  "
Set-TextValue $ws.Range("D35") "149.46"
$ws.Range("E35").Value = "  -2.07%  "
Set-TextValue $ws.Range("D36") "4.04"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("E37").Value = "  -2.72%  "
Set-TextValue $ws.Range("D38") "0.876"
$ws.Range("E38").Value = "  -3.08%  "
Set-TextValue $ws.Range("D39") "0.839"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D40") "1.46"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D41") "35.92"
$ws.Range("E41").Value = "  -2.29%  "
Set-TextValue $ws.Range("D42") "3.55"
$ws.Range("E42").Value = "  -2.00%  "
Set-TextValue $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  -0.01%  "
Set-TextValue $ws.Range("D44") "274.68"
$ws.Range("E44").Value = "  +2.21%  "
$ws.Range("E45").Value = "  -2.43%  "
Set-TextValue $ws.Range("D46") "10.66"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("D50").Value = "1.974.32"
$ws.Range("E50").Value = "  -3.30%  "
Set-TextValue $ws.Range("D51") "18.56"
$ws.Range("E51").Value = "  +1.65%  "
